# Data Driven test data update for trello testdata.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "Boardname"
$ws.Range("B1").Value = "CardTitle"
$ws.Range("C1").Value = "CardDescription"

# --- Data rows: Boardname | CardTitle | CardDescription ---
$ws.Range("A2").Value = "Board_1453"
$ws.Range("B2").Value = "Learning Java"
$ws.Range("C2").Value = "Learning Java Description"

$ws.Range("A3").Value = "Board_1453"
$ws.Range("B3").Value = "Learning TestNG"
$ws.Range("C3").Value = "Learning TestNG Description"

$ws.Range("A4").Value = "Board_1453"
$ws.Range("B4").Value = "Learning RestAssured"
$ws.Range("C4").Value = "Learning RestAssured Description"

$ws.Range("A5").Value = "Board_1453"
$ws.Range("B5").Value = "Learning Cucumber"
$ws.Range("C5").Value = "Learning Cucumber Description"

$ws.Range("A6").Value = "Board_1453"
$ws.Range("B6").Value = "Learning Mockito"
$ws.Range("C6").Value = "Learning Mockito Description"

$ws.Range("A7").Value = "Board_1453"
$ws.Range("B7").Value = "Learning XYZ"
$ws.Range("C7").Value = "Learning XYZ Description"

# --- Column widths (stored width = ColumnWidth + 0.8333333) ---
$ws.Columns.Item(2).ColumnWidth = 18.1666666666667
$ws.Columns.Item(3).ColumnWidth = 28.1666666666667

# --- Selection moves from A2:A7 to B2:B7 ---
$ws.Range("B2:B7").Select() | Out-Null
